$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new BOM row (row 14) for the slide switch power switch.
# Copy formatting from the row above first so styles (currency number
# format on Price/Total, hyperlink style on Link) carry over cleanly.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

# Fill in the new row's values. Order matters for shared-string layout:
# Link text, then Description, then Part (matches how the source file
# appends new shared strings).
$ws.Range("F14").Value = "MK-12D18-G020 G-Switch | C3019727 - LCSC Electronics"
$ws.Range("B14").Value = "Slide switch"
$ws.Range("A14").Value = "MK-12D18-G020"
$ws.Range("C14").Value = 0.63
$ws.Range("D14").Value = 1
$ws.Range("E14").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

# Hook up the LCSC product-page hyperlink for the new Link cell.
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lcsc.com/product-detail/Slide-Switches_G-Switch-MK-12D18-G020_C3019727.html", [Type]::Missing, [Type]::Missing, "MK-12D18-G020 G-Switch | C3019727 - LCSC Electronics")

# Grow Table1 so the new row becomes part of the structured table.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F14"))

# Reflect the in-progress routing work: selection now sits at B21:B22.
$ws.Range("B21:B22").Select()
